$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 209
$ws1.Range("F4").Value = 652
$ws1.Range("F5").Value = 567
$ws1.Range("F6").Value = 314
$ws1.Range("F7").Value = 2791
$ws1.Range("F8").Value = 470
$ws1.Range("F9").Value = 7865
$ws1.Range("F10").Value = 204
$ws1.Range("F11").Value = 472
$ws1.Range("F13").Value = 361
$ws1.Range("F14").Value = 49

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 209
$ws4.Range("F4").Value = 652
$ws4.Range("F5").Value = 567
$ws4.Range("F6").Value = 314
$ws4.Range("F9").Value = 2791
$ws4.Range("F10").Value = 470
$ws4.Range("F11").Value = 7865
$ws4.Range("F12").Value = 204
$ws4.Range("F13").Value = 472
$ws4.Range("F17").Value = 361
$ws4.Range("F18").Value = 49
